$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 961.2
$ws.Range("I43").Value = 825
$ws.Range("J43").Value = 1052
$ws.Range("K43").Value = 825
$ws.Range("L43").Value = 1052
$ws.Range("M43").Value = -756
$ws.Range("N43").Value = -1190

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2507.1428
$ws.Range("I51").Value = 2237.5
$ws.Range("J51").Value = 2866.6667
$ws.Range("K51").Value = 2237.5
$ws.Range("L51").Value = 2866.6667
$ws.Range("M51").Value = -1753.5
$ws.Range("N51").Value = -3834.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3637.8816
$ws.Range("I64").Value = 3390.476
$ws.Range("J64").Value = 3943.5
$ws.Range("K64").Value = 3390.476
$ws.Range("L64").Value = 3943.5
$ws.Range("M64").Value = -3142.476
$ws.Range("N64").Value = -4439.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3637.8816
$ws.Range("I67").Value = 3390.476
$ws.Range("J67").Value = 3943.5
$ws.Range("K67").Value = 3390.476
$ws.Range("L67").Value = 3943.5
$ws.Range("M67").Value = -2532.476
$ws.Range("N67").Value = -5659.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3109.9473
$ws.Range("I76").Value = 2338.125
$ws.Range("J76").Value = 3671.2727
$ws.Range("K76").Value = 2338.125
$ws.Range("L76").Value = 3671.2727
$ws.Range("M76").Value = -2023.125
$ws.Range("N76").Value = -4301.2727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3109.9473
$ws.Range("I79").Value = 2338.125
$ws.Range("J79").Value = 3671.2727
$ws.Range("K79").Value = 2338.125
$ws.Range("L79").Value = 3671.2727
$ws.Range("M79").Value = -1246.125
$ws.Range("N79").Value = -5855.2727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 207.2
$ws.Range("I99").Value = 207.2
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 621.5999999999999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 876.4000000000001
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16683662
$ws.Range("I32").Value = 25005918
$ws.Range("J32").Value = 39149.95
$ws.Range("K32").Value = 25005918
$ws.Range("L32").Value = 39149.95
$ws.Range("M32").Value = -25005631
$ws.Range("N32").Value = -39723.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2067.182
$ws.Range("I61").Value = 1514.3928
$ws.Range("J61").Value = 5162.8
$ws.Range("K61").Value = 1514.3928
$ws.Range("L61").Value = 5162.8
$ws.Range("M61").Value = -1302.3928
$ws.Range("N61").Value = -5586.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2067.182
$ws.Range("I136").Value = 1514.3928
$ws.Range("J136").Value = 5162.8
$ws.Range("K136").Value = 4543.178400000001
$ws.Range("L136").Value = 15488.4
$ws.Range("M136").Value = -1993.178400000001
$ws.Range("N136").Value = -20588.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1670.8572
$ws.Range("I99").Value = 1039.8235
$ws.Range("J99").Value = 4352.75
$ws.Range("K99").Value = 1039.8235
$ws.Range("L99").Value = 4352.75
$ws.Range("M99").Value = 458.1765
$ws.Range("N99").Value = -7348.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 28509
$ws.Range("I64").Value = 15256
$ws.Range("J64").Value = 35135.5
$ws.Range("K64").Value = 15256
$ws.Range("L64").Value = 35135.5
$ws.Range("M64").Value = -15008
$ws.Range("N64").Value = -35631.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 28509
$ws.Range("I67").Value = 15256
$ws.Range("J67").Value = 35135.5
$ws.Range("K67").Value = 15256
$ws.Range("L67").Value = 35135.5
$ws.Range("M67").Value = -14398
$ws.Range("N67").Value = -36851.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 753.0769
$ws.Range("I105").Value = 679
$ws.Range("K105").Value = 679
$ws.Range("M105").Value = 1068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 713.1042
$ws.Range("I131").Value = 408.44446
$ws.Range("J131").Value = 895.9
$ws.Range("K131").Value = 1225.33338
$ws.Range("L131").Value = 2687.7
$ws.Range("M131").Value = 3814.66662
$ws.Range("N131").Value = -12767.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1012
$ws.Range("I23").Value = 1012
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1012
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -789
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 980
$ws.Range("I27").Value = 980
$ws.Range("K27").Value = 980
$ws.Range("M27").Value = -814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3076.72
$ws.Range("I80").Value = 2825
$ws.Range("J80").Value = 3524.2222
$ws.Range("K80").Value = 2825
$ws.Range("L80").Value = 3524.2222
$ws.Range("M80").Value = -1827
$ws.Range("N80").Value = -5520.2222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3076.72
$ws.Range("I83").Value = 2825
$ws.Range("J83").Value = 3524.2222
$ws.Range("K83").Value = 14125
$ws.Range("L83").Value = 17621.111
$ws.Range("M83").Value = -9133
$ws.Range("N83").Value = -27605.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3992.3076
$ws.Range("J2").Value = 3992.3076
$ws.Range("L2").Value = 3992.3076
$ws.Range("N2").Value = -4216.3076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2688.625
$ws.Range("I68").Value = 2813.1765
$ws.Range("J68").Value = 2386.1428
$ws.Range("K68").Value = 2813.1765
$ws.Range("L68").Value = 2386.1428
$ws.Range("M68").Value = -2064.1765
$ws.Range("N68").Value = -3884.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2688.625
$ws.Range("I71").Value = 2813.1765
$ws.Range("J71").Value = 2386.1428
$ws.Range("K71").Value = 14065.8825
$ws.Range("L71").Value = 11930.714
$ws.Range("M71").Value = -10321.8825
$ws.Range("N71").Value = -19418.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2990.875
$ws.Range("I82").Value = 2360.4
$ws.Range("J82").Value = 4041.6667
$ws.Range("K82").Value = 2360.4
$ws.Range("L82").Value = 4041.6667
$ws.Range("M82").Value = -1999.4
$ws.Range("N82").Value = -4763.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2990.875
$ws.Range("I85").Value = 2360.4
$ws.Range("J85").Value = 4041.6667
$ws.Range("K85").Value = 2360.4
$ws.Range("L85").Value = 4041.6667
$ws.Range("M85").Value = -1112.4
$ws.Range("N85").Value = -6537.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2324.7
$ws.Range("I132").Value = 1895.7142
$ws.Range("J132").Value = 3760.8696
$ws.Range("K132").Value = 5687.142599999999
$ws.Range("L132").Value = 11282.6088
$ws.Range("M132").Value = -3157.142599999999
$ws.Range("N132").Value = -16342.6088

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12514612
$ws.Range("I62").Value = 16684816
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 16684816
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -16684192
$ws.Range("N62").Value = -5248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 20249
$ws.Range("J63").Value = 20249
$ws.Range("L63").Value = 20249
$ws.Range("N63").Value = -21497

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 12514612
$ws.Range("I65").Value = 16684816
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 83424080
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -83420960
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 20249
$ws.Range("J66").Value = 20249
$ws.Range("L66").Value = 60747
$ws.Range("N66").Value = -66987
